# Apply updated odds values to Sheet1, as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 6.8
$ws.Range("J2").Value = 14
$ws.Range("K2").Value = 2.7
$ws.Range("T2").Value = 3.28
$ws.Range("U2").Value = 2.9
$ws.Range("V2").Value = 1.31
$ws.Range("AB2").Value = 400
$ws.Range("AC2").Value = 13
$ws.Range("AG2").Value = 6.5
$ws.Range("AH2").Value = 4.9
$ws.Range("AK2").Value = 13.5
$ws.Range("AT2").Value = 2.72
$ws.Range("AV2").Value = 22

# Row 3
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 3.5
$ws.Range("I3").Value = 2.05
$ws.Range("J3").Value = 3.6
$ws.Range("L3").Value = 2.75
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 11
$ws.Range("O3").Value = 1.21
$ws.Range("P3").Value = 3.5
$ws.Range("Q3").Value = 1.67
$ws.Range("R3").Value = 2.15
$ws.Range("S3").Value = 1.32
$ws.Range("T3").Value = 3.2
$ws.Range("U3").Value = 1.55
$ws.Range("V3").Value = 2.15
$ws.Range("W3").Value = 12
$ws.Range("X3").Value = 17
$ws.Range("Y3").Value = 11
$ws.Range("Z3").Value = 34
$ws.Range("AA3").Value = 23
$ws.Range("AB3").Value = 29
$ws.Range("AG3").Value = 9.5
$ws.Range("AH3").Value = 11
$ws.Range("AI3").Value = 9
$ws.Range("AJ3").Value = 21
$ws.Range("AK3").Value = 17
$ws.Range("AL3").Value = 23
$ws.Range("AN3").Value = 5.5
$ws.Range("AO3").Value = 17
$ws.Range("AP3").Value = 23
$ws.Range("AQ3").Value = 3.25
$ws.Range("AS3").Value = 41
$ws.Range("AT3").Value = 4.33
$ws.Range("AU3").Value = 11
$ws.Range("AV3").Value = 19
$ws.Range("AZ3").Value = 51
$ws.Range("BA3").Value = 67
$ws.Range("BB3").Value = 151

# Row 4
$ws.Range("G4").Value = 3.6
$ws.Range("J4").Value = 4.5
$ws.Range("V4").Value = 1.73
$ws.Range("X4").Value = 17
$ws.Range("AC4").Value = 9
$ws.Range("AG4").Value = 6.5
$ws.Range("BA4").Value = 101

# Row 5
$ws.Range("G5").Value = 2.18
$ws.Range("H5").Value = 2.75
$ws.Range("I5").Value = 3.9
$ws.Range("J5").Value = 2.82
$ws.Range("K5").Value = 1.88
$ws.Range("M5").Value = 1.13
$ws.Range("N5").Value = 5.2
$ws.Range("O5").Value = 1.52
$ws.Range("P5").Value = 2.37
$ws.Range("Q5").Value = 2.52
$ws.Range("R5").Value = 1.47
$ws.Range("S5").Value = 1.55
$ws.Range("T5").Value = 2.3
$ws.Range("U5").Value = 2.05
$ws.Range("V5").Value = 1.7
$ws.Range("W5").Value = 5.7
$ws.Range("X5").Value = 9.25
$ws.Range("Y5").Value = 9
$ws.Range("Z5").Value = 22
$ws.Range("AA5").Value = 21
$ws.Range("AB5").Value = 40
$ws.Range("AC5").Value = 5.2
$ws.Range("AD5").Value = 5.5
$ws.Range("AE5").Value = 17
$ws.Range("AF5").Value = 110
$ws.Range("AG5").Value = 8.25
$ws.Range("AH5").Value = 19.5
$ws.Range("AJ5").Value = 65
$ws.Range("AL5").Value = 60
$ws.Range("AN5").Value = 3.9
$ws.Range("AO5").Value = 11.75
$ws.Range("AP5").Value = 22
$ws.Range("AQ5").Value = 2.3
$ws.Range("AR5").Value = 7.4
$ws.Range("AS5").Value = 80
$ws.Range("AT5").Value = 5.6
$ws.Range("AV5").Value = 32
$ws.Range("AX5").Value = 200
$ws.Range("AY5").Value = 500
$ws.Range("AZ5").Value = 50
$ws.Range("BA5").Value = 100
$ws.Range("BB5").Value = 350

# Row 7
$ws.Range("G7").Value = 2.02
$ws.Range("H7").Value = 3.2
$ws.Range("I7").Value = 3.7
$ws.Range("J7").Value = 2.7
$ws.Range("L7").Value = 4.05
$ws.Range("N7").Value = 7.3
$ws.Range("O7").Value = 1.28
$ws.Range("P7").Value = 3.35
$ws.Range("Q7").Value = 1.83
$ws.Range("R7").Value = 1.87
$ws.Range("S7").Value = 1.42
$ws.Range("T7").Value = 2.65
$ws.Range("U7").Value = 1.65
$ws.Range("V7").Value = 2.1
$ws.Range("W7").Value = 7.5
$ws.Range("X7").Value = 9.75
$ws.Range("Z7").Value = 19
$ws.Range("AA7").Value = 16
$ws.Range("AB7").Value = 25
$ws.Range("AC7").Value = 7.3
$ws.Range("AD7").Value = 6.2
$ws.Range("AE7").Value = 12.5
$ws.Range("AF7").Value = 50
$ws.Range("AG7").Value = 12
$ws.Range("AH7").Value = 22
$ws.Range("AI7").Value = 11.75
$ws.Range("AK7").Value = 30
$ws.Range("AL7").Value = 32
$ws.Range("AM7").Value = 350
$ws.Range("AN7").Value = 3.95
$ws.Range("AO7").Value = 11
$ws.Range("AP7").Value = 19.5
$ws.Range("AQ7").Value = 2.65
$ws.Range("AR7").Value = 6.8
$ws.Range("AS7").Value = 60
$ws.Range("AT7").Value = 5.5
$ws.Range("AU7").Value = 20
$ws.Range("AV7").Value = 24
$ws.Range("AW7").Value = 100
$ws.Range("AX7").Value = 120
$ws.Range("AY7").Value = 300
$ws.Range("AZ7").Value = 45
$ws.Range("BA7").Value = 80

# Row 8
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 8
$ws.Range("Q8").Value = 2.15
$ws.Range("R8").Value = 1.67
$ws.Range("W8").Value = 5
$ws.Range("AC8").Value = 8
$ws.Range("AF8").Value = 126
$ws.Range("AJ8").Value = 126
$ws.Range("AK8").Value = 81
$ws.Range("AL8").Value = 81
$ws.Range("AP8").Value = 23
$ws.Range("BB8").Value = 251

# Row 9
$ws.Range("L9").Value = 4
$ws.Range("N9").Value = 7.5
$ws.Range("Z9").Value = 21
$ws.Range("AB9").Value = 34

# Row 10
$ws.Range("G10").Value = 1.57
$ws.Range("H10").Value = 3.4
$ws.Range("I10").Value = 5.25
$ws.Range("J10").Value = 2.3
$ws.Range("K10").Value = 2.05
$ws.Range("L10").Value = 6.5
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 7.5
$ws.Range("Q10").Value = 2.35
$ws.Range("R10").Value = 1.57
$ws.Range("S10").Value = 1.5
$ws.Range("T10").Value = 2.5
$ws.Range("U10").Value = 2.38
$ws.Range("V10").Value = 1.53
$ws.Range("W10").Value = 5
$ws.Range("X10").Value = 6.5
$ws.Range("Z10").Value = 11
$ws.Range("AC10").Value = 7.5
$ws.Range("AD10").Value = 7.5
$ws.Range("AE10").Value = 23
$ws.Range("AG10").Value = 12
$ws.Range("AH10").Value = 29
$ws.Range("AI10").Value = 21
$ws.Range("AJ10").Value = 67
$ws.Range("AN10").Value = 3.4
$ws.Range("AO10").Value = 8.5
$ws.Range("AQ10").Value = 2.5
$ws.Range("AT10").Value = 7.5
$ws.Range("AW10").Value = 151
$ws.Range("AX10").Value = 201
$ws.Range("AZ10").Value = 29
$ws.Range("BA10").Value = 51

# Row 12
$ws.Range("M12").Value = 1.07
$ws.Range("N12").Value = 9
$ws.Range("O12").Value = 1.36
$ws.Range("P12").Value = 3
$ws.Range("Q12").Value = 2.2
$ws.Range("R12").Value = 1.65

# Row 13
$ws.Range("O13").Value = 1.25
$ws.Range("P13").Value = 3.75
$ws.Range("Q13").Value = 1.9
$ws.Range("R13").Value = 1.95

# Row 14
$ws.Range("H14").Value = 3.7
$ws.Range("I14").Value = 3.3
$ws.Range("M14").Value = 1.03
$ws.Range("N14").Value = 15
$ws.Range("Q14").Value = 1.65
$ws.Range("R14").Value = 2.2
$ws.Range("S14").Value = 1.3
$ws.Range("T14").Value = 3.4
$ws.Range("W14").Value = 10
$ws.Range("AC14").Value = 15
$ws.Range("AD14").Value = 7.5
$ws.Range("AK14").Value = 23
$ws.Range("AM14").Value = 126
$ws.Range("AQ14").Value = 3.4
$ws.Range("AY14").Value = 126
$ws.Range("BB14").Value = 101

# Row 15
$ws.Range("Q15").Value = 3.1
$ws.Range("R15").Value = 1.36

